# Update Koffi benchmark results (doc/benchmarks.xlsx)
#
# The benchmark numbers moved for a few Koffi / node-ffi data points on
# both the Linux and Windows sheets. Update the raw "Time" column (C);
# the "D" column display formula (=_xlfn.CONCAT("x ", SUBSTITUTE(...)))
# recalculates on its own.

$wb = $excel.ActiveWorkbook

$wsLinux = $wb.Worksheets.Item("Linux")
$wsWindows = $wb.Worksheets.Item("Windows")

# --- Linux sheet: atoi/Koffi, atoi/node-ffi, raylib/node-ffi ---
$wsLinux.Range("C8").Value = 0.62
$wsLinux.Range("C9").Value = 0.009
$wsLinux.Range("C13").Value = 0.28

# --- Windows sheet: rand/Koffi, atoi/Koffi, atoi/node-ffi, raylib/Koffi ---
$wsWindows.Range("C4").Value = 0.77
$wsWindows.Range("C8").Value = 0.62
$wsWindows.Range("C9").Value = 0.009
$wsWindows.Range("C12").Value = 0.83

# --- Selection / active sheet bookkeeping (matches the saved view state) ---
$wsLinux.Range("C20").Select() | Out-Null
$wsWindows.Range("N7").Select() | Out-Null
$wsWindows.Activate() | Out-Null
